$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$newShp = $s.Shapes.AddTextbox(1, 100, 100, 300, 100)
$tr = $newShp.TextFrame.TextRange
$tr.Text = ""
Write-Host "empty text len:" $tr.Text.Length
$tr.Font.Size = 33
Write-Host "Font size now:" $tr.Font.Size
